$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (MuSCs target rows removed, table shrinks from A1:T10 to A1:T7)
$ws.Range("A8:T10").Delete()

# Update remaining rows 2-7 with the new TPM-derived values
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cadm1"
$ws.Range("C2").Value = "Crtam"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.019613333333333
$ws.Range("H2").Value = 12.05884
$ws.Range("I2").Value = 0.4377217086785624
$ws.Range("J2").Value = 0.4377217086785624
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1197753333333333
$ws.Range("N2").Value = 0.359326
$ws.Range("O2").Value = 0.9004197311214965
$ws.Range("P2").Value = 0.9004197311214965
$ws.Range("Q2").Value = 0.4814505268711111
$ws.Range("R2").Value = 4.33305474184
$ws.Range("S2").Value = 0.3941332632343932
$ws.Range("T2").Value = 0.3941332632343932
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cadm1"
$ws.Range("C3").Value = "Crtam"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.019613333333333
$ws.Range("H3").Value = 12.05884
$ws.Range("I3").Value = 0.4377217086785624
$ws.Range("J3").Value = 0.4377217086785624
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01324633333333333
$ws.Range("N3").Value = 0.039739
$ws.Range("O3").Value = 0.09958026887850349
$ws.Range("P3").Value = 0.09958026887850349
$ws.Range("Q3").Value = 0.05324513808444443
$ws.Range("R3").Value = 0.4792062427599999
$ws.Range("S3").Value = 0.04358844544416921
$ws.Range("T3").Value = 0.04358844544416921
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cadm1"
$ws.Range("C4").Value = "Crtam"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.236474
$ws.Range("H4").Value = 0.709422
$ws.Range("I4").Value = 0.02575118419467902
$ws.Range("J4").Value = 0.02575118419467902
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1197753333333333
$ws.Range("N4").Value = 0.359326
$ws.Range("O4").Value = 0.9004197311214965
$ws.Range("P4").Value = 0.9004197311214965
$ws.Range("Q4").Value = 0.02832375217466667
$ws.Range("R4").Value = 0.254913769572
$ws.Range("S4").Value = 0.02318687434863301
$ws.Range("T4").Value = 0.02318687434863301
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cadm1"
$ws.Range("C5").Value = "Crtam"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.236474
$ws.Range("H5").Value = 0.709422
$ws.Range("I5").Value = 0.02575118419467902
$ws.Range("J5").Value = 0.02575118419467902
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01324633333333333
$ws.Range("N5").Value = 0.039739
$ws.Range("O5").Value = 0.09958026887850349
$ws.Range("P5").Value = 0.09958026887850349
$ws.Range("Q5").Value = 0.003132413428666666
$ws.Range("R5").Value = 0.028191720858
$ws.Range("S5").Value = 0.002564309846046006
$ws.Range("T5").Value = 0.002564309846046006
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Cadm1"
$ws.Range("C6").Value = "Crtam"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.926946666666667
$ws.Range("H6").Value = 14.78084
$ws.Range("I6").Value = 0.5365271071267587
$ws.Range("J6").Value = 0.5365271071267587
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1197753333333333
$ws.Range("N6").Value = 0.359326
$ws.Range("O6").Value = 0.9004197311214965
$ws.Range("P6").Value = 0.9004197311214965
$ws.Range("Q6").Value = 0.5901266793155556
$ws.Range("R6").Value = 5.311140113840001
$ws.Range("S6").Value = 0.4830995935384704
$ws.Range("T6").Value = 0.4830995935384704
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cadm1"
$ws.Range("C7").Value = "Crtam"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.926946666666667
$ws.Range("H7").Value = 14.78084
$ws.Range("I7").Value = 0.5365271071267587
$ws.Range("J7").Value = 0.5365271071267587
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01324633333333333
$ws.Range("N7").Value = 0.039739
$ws.Range("O7").Value = 0.09958026887850349
$ws.Range("P7").Value = 0.09958026887850349
$ws.Range("Q7").Value = 0.06526397786222222
$ws.Range("R7").Value = 0.58737580076
$ws.Range("S7").Value = 0.05342751358828828
$ws.Range("T7").Value = 0.05342751358828828
